$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (cluster 1 head)
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 0

# Row 4
$ws.Range("A4").Value = 1
$ws.Range("C4").Value = 0

# Row 5 (cluster 2 head)
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 1

# Row 6
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 0

# Row 7
$ws.Range("D7").Value = 1
$ws.Range("F7").Value = 0

# Row 8 (cluster 3 head)
$ws.Range("H8").Value = 1

# Row 9
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 0

# Row 10
$ws.Range("I10").Value = 0
